$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "G2" = 118.0346986666667
    "H2" = 354.104096
    "I2" = 0.2666057129183408
    "J2" = 0.2666057129183408
    "M2" = 19.84999933333333
    "N2" = 59.549998
    "O2" = 0.1831667009459596
    "P2" = 0.1831667009459596
    "Q2" = 2342.988689843534
    "R2" = 21086.89820859181
    "S2" = 0.04883328888859809
    "T2" = 0.04883328888859809
    "G3" = 118.0346986666667
    "H3" = 354.104096
    "I3" = 0.2666057129183408
    "J3" = 0.2666057129183408
    "O3" = 0.3072686534975208
    "P3" = 0.3072686534975207
    "Q3" = 3930.446834332327
    "R3" = 35374.02150899095
    "S3" = 0.08191957842316518
    "T3" = 0.08191957842316516
    "G4" = 118.0346986666667
    "H4" = 354.104096
    "I4" = 0.2666057129183408
    "J4" = 0.2666057129183408
    "M4" = 24.07451633333333
    "N4" = 72.22354899999999
    "O4" = 0.2221486086521591
    "P4" = 0.2221486086521591
    "Q4" = 2841.628280950744
    "R4" = 25574.6545285567
    "S4" = 0.05922608818352638
    "T4" = 0.05922608818352638
    "G5" = 118.0346986666667
    "H5" = 354.104096
    "I5" = 0.2666057129183408
    "J5" = 0.2666057129183408
    "M5" = 31.14762733333333
    "N5" = 93.442882
    "O5" = 0.2874160369043605
    "P5" = 0.2874160369043605
    "Q5" = 3676.50080647163
    "R5" = 33088.50725824467
    "S5" = 0.07662675742305118
    "T5" = 0.07662675742305118
    "I6" = 0.4881754016778185
    "J6" = 0.4881754016778186
    "M6" = 19.84999933333333
    "N6" = 59.549998
    "O6" = 0.1831667009459596
    "P6" = 0.1831667009459596
    "Q6" = 4290.191055062974
    "R6" = 38611.71949556677
    "S6" = 0.08941747780829469
    "T6" = 0.0894174778082947
    "I7" = 0.4881754016778185
    "J7" = 0.4881754016778186
    "O7" = 0.3072686534975208
    "P7" = 0.3072686534975207
    "S7" = 0.1500009983441546
    "T7" = 0.1500009983441546
    "I8" = 0.4881754016778185
    "J8" = 0.4881754016778186
    "M8" = 24.07451633333333
    "N8" = 72.22354899999999
    "O8" = 0.2221486086521591
    "P8" = 0.2221486086521591
    "Q8" = 5203.238191287636
    "R8" = 46829.14372158873
    "S8" = 0.1084474862609363
    "T8" = 0.1084474862609363
    "I9" = 0.4881754016778185
    "J9" = 0.4881754016778186
    "M9" = 31.14762733333333
    "N9" = 93.442882
    "O9" = 0.2874160369043605
    "P9" = 0.2874160369043605
    "Q9" = 6731.953484124466
    "R9" = 60587.58135712019
    "S9" = 0.1403094392644329
    "T9" = 0.1403094392644329
    "G10" = 45.876452
    "H10" = 137.629356
    "I10" = 0.1036214293744632
    "J10" = 0.1036214293744632
    "M10" = 19.84999933333333
    "N10" = 59.549998
    "O10" = 0.1831667009459596
    "P10" = 0.1831667009459596
    "Q10" = 910.6475416156986
    "R10" = 8195.827874541288
    "S10" = 0.01897999536582517
    "T10" = 0.01897999536582518
    "G11" = 45.876452
    "H11" = 137.629356
    "I11" = 0.1036214293744632
    "J11" = 0.1036214293744632
    "O11" = 0.3072686534975208
    "P11" = 0.3072686534975207
    "Q11" = 1527.64362997201
    "R11" = 13748.79266974809
    "S11" = 0.03183961707737975
    "T11" = 0.03183961707737975
    "G12" = 45.876452
    "H12" = 137.629356
    "I12" = 0.1036214293744632
    "J12" = 0.1036214293744632
    "M12" = 24.07451633333333
    "N12" = 72.22354899999999
    "O12" = 0.2221486086521591
    "P12" = 0.2221486086521591
    "Q12" = 1104.453392989383
    "R12" = 9940.080536904443
    "S12" = 0.02301935636208497
    "T12" = 0.02301935636208497
    "G13" = 45.876452
    "H13" = 137.629356
    "I13" = 0.1036214293744632
    "J13" = 0.1036214293744632
    "M13" = 31.14762733333333
    "N13" = 93.442882
    "O13" = 0.2874160369043605
    "P13" = 0.2874160369043605
    "Q13" = 1428.942630271555
    "R13" = 12860.48367244399
    "S13" = 0.02978246056917329
    "T13" = 0.0297824605691733
    "G14" = 62.68962833333333
    "H14" = 188.068885
    "I14" = 0.1415974560293775
    "J14" = 0.1415974560293775
    "M14" = 19.84999933333333
    "N14" = 59.549998
    "O14" = 0.1831667009459596
    "P14" = 0.1831667009459596
    "Q14" = 1244.389080623581
    "R14" = 11199.50172561223
    "S14" = 0.02593593888324165
    "T14" = 0.02593593888324165
    "G15" = 62.68962833333333
    "H15" = 188.068885
    "I15" = 0.1415974560293775
    "J15" = 0.1415974560293775
    "O15" = 0.3072686534975208
    "P15" = 0.3072686534975207
    "Q15" = 2087.506928145391
    "R15" = 18787.56235330852
    "S15" = 0.04350845965282122
    "T15" = 0.04350845965282122
    "G16" = 62.68962833333333
    "H16" = 188.068885
    "I16" = 0.1415974560293775
    "J16" = 0.1415974560293775
    "M16" = 24.07451633333333
    "N16" = 72.22354899999999
    "O16" = 0.2221486086521591
    "P16" = 0.2221486086521591
    "Q16" = 1509.222481241429
    "R16" = 13583.00233117286
    "S16" = 0.03145567784561148
    "T16" = 0.03145567784561149
    "G17" = 62.68962833333333
    "H17" = 188.068885
    "I17" = 0.1415974560293775
    "J17" = 0.1415974560293775
    "M17" = 31.14762733333333
    "N17" = 93.442882
    "O17" = 0.2874160369043605
    "P17" = 0.2874160369043605
    "Q17" = 1952.633180991841
    "R17" = 17573.69862892657
    "S17" = 0.04069737964770311
    "T17" = 0.04069737964770312
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
